{"js": "const pairs = [\n  [\"2023-04-12 Wednesday\", \"2023-04-13 Thursday\"],\n  [\"4+37=\", \"77-67=\"],\n  [\"32-27=\", \"23+36=\"],\n  [\"75-30=\", \"59+24=\"],\n  [\"44+39=\", \"68+28=\"],\n  [\"32+42=\", \"36+7=\"],\n  [\"98-83=\", \"32+12=\"],\n  [\"88-49=\", \"0+81=\"],\n  [\"83-3=\", \"67-55=\"],\n  [\"59-27=\", \"43-7=\"],\n  [\"97-22=\", \"4+93=\"],\n  [\"53-32=\", \"23+11=\"],\n  [\"11+67=\", \"61+28=\"],\n  [\"10+81=\", \"33+39=\"],\n  [\"75-31=\", \"24+31=\"],\n  [\"28+29=\", \"63+13=\"],\n  [\"99-88=\", \"60-30=\"],\n  [\"93-69=\", \"87-23=\"],\n  [\"46+36=\", \"23+48=\"],\n  [\"37+61=\", \"80-71=\"],\n  [\"60-41=\", \"6+36=\"],\n  [\"77-54=\", \"6+12=\"],\n  [\"40-29=\", \"86-62=\"],\n  [\"89-5=\", \"2+68=\"],\n  [\"85-53=\", \"37+57=\"],\n  [\"96-71=\", \"59+20=\"],\n  [\"96-59=\", \"54-50=\"],\n  [\"10+63=\", \"44+50=\"],\n  [\"87-44=\", \"41+3=\"],\n  [\"66+10=\", \"44+46=\"],\n  [\"6+34=\", \"68+26=\"],\n  [\"12+34=\", \"9+6=\"],\n  [\"60+23=\", \"60-0=\"],\n  [\"23+68=\", \"48-32=\"],\n  [\"28+27=\", \"39-28=\"],\n  [\"51+19=\", \"87-48=\"],\n  [\"58+9=\", \"54+8=\"],\n  [\"4+47=\", \"20-11=\"],\n  [\"46+14=\", \"41+22=\"],\n  [\"15-11=\", \"54+11=\"],\n  [\"34-5=\", \"50-9=\"],\n  [\"42+11=\", \"7+83=\"],\n  [\"10+61=\", \"99-90=\"],\n  [\"98-76=\", \"20+2=\"],\n  [\"95-78=\", \"32+22=\"],\n  [\"38-19=\", \"95-55=\"],\n  [\"68-60=\", \"7+26=\"],\n  [\"81-64=\", \"34-6=\"],\n  [\"26+9=\", \"26+12=\"],\n  [\"71-49=\", \"90-7=\"],\n  [\"17+26=\", \"11+9=\"],\n  [\"44+14=\", \"97-31=\"],\n  [\"82-3=\", \"57+41=\"],\n  [\"53+38=\", \"59+4=\"],\n  [\"27-27=\", \"2+78=\"],\n  [\"45+51=\", \"47+35=\"],\n  [\"31+4=\", \"42+18=\"],\n  [\"48-26=\", \"54+10=\"],\n  [\"62-13=\", \"70-54=\"],\n  [\"67-16=\", \"93-55=\"],\n  [\"10-3=\", \"57-12=\"],\n  [\"51+33=\", \"16+59=\"],\n  [\"91-34=\", \"70+10=\"],\n  [\"30+69=\", \"88-7=\"],\n  [\"48+3=\", \"54-21=\"],\n  [\"73-41=\", \"62+23=\"],\n  [\"6+86=\", \"35+41=\"],\n  [\"6+48=\", \"1+17=\"],\n  [\"79-78=\", \"64-28=\"],\n  [\"54+29=\", \"74+8=\"],\n  [\"25+14=\", \"35+52=\"],\n  [\"7+70=\", \"7+34=\"],\n  [\"55-22=\", \"42-37=\"],\n  [\"27+37=\", \"24+3=\"],\n  [\"70-18=\", \"87-77=\"],\n  [\"87-76=\", \"27+38=\"],\n  [\"70+22=\", \"64+10=\"],\n  [\"93-51=\", \"37-12=\"],\n  [\"24+38=\", \"23+11=\"],\n  [\"48-13=\", \"97-54=\"],\n  [\"10+33=\", \"21+11=\"],\n  [\"87-3=\", \"9+44=\"],\n  [\"92-90=\", \"94-67=\"],\n  [\"24+4=\", \"43+36=\"],\n  [\"23-1=\", \"31+46=\"],\n  [\"88-60=\", \"52+18=\"],\n  [\"64-18=\", \"14+26=\"],\n  [\"80+5=\", \"6+93=\"],\n  [\"9+16=\", \"19+80=\"],\n  [\"90-53=\", \"24+17=\"],\n  [\"98-90=\", \"2+30=\"],\n  [\"88-17=\", \"41-22=\"],\n  [\"14+17=\", \"15+31=\"],\n  [\"0+83=\", \"51-31=\"],\n  [\"29+36=\", \"73-69=\"],\n  [\"79-48=\", \"24-17=\"],\n  [\"63-9=\", \"34-17=\"],\n  [\"97-66=\", \"34+61=\"],\n  [\"92-60=\", \"97-38=\"],\n  [\"23+26=\", \"72+1=\"],\n  [\"2+72=\", \"52-1=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Not found: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, 'Replace');\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@('2023-04-12 Wednesday', '2023-04-13 Thursday')\n  ,@('4+37=', '77-67=')\n  ,@('32-27=', '23+36=')\n  ,@('75-30=', '59+24=')\n  ,@('44+39=', '68+28=')\n  ,@('32+42=', '36+7=')\n  ,@('98-83=', '32+12=')\n  ,@('88-49=', '0+81=')\n  ,@('83-3=', '67-55=')\n  ,@('59-27=', '43-7=')\n  ,@('97-22=', '4+93=')\n  ,@('53-32=', '23+11=')\n  ,@('11+67=', '61+28=')\n  ,@('10+81=', '33+39=')\n  ,@('75-31=', '24+31=')\n  ,@('28+29=', '63+13=')\n  ,@('99-88=', '60-30=')\n  ,@('93-69=', '87-23=')\n  ,@('46+36=', '23+48=')\n  ,@('37+61=', '80-71=')\n  ,@('60-41=', '6+36=')\n  ,@('77-54=', '6+12=')\n  ,@('40-29=', '86-62=')\n  ,@('89-5=', '2+68=')\n  ,@('85-53=', '37+57=')\n  ,@('96-71=', '59+20=')\n  ,@('96-59=', '54-50=')\n  ,@('10+63=', '44+50=')\n  ,@('87-44=', '41+3=')\n  ,@('66+10=', '44+46=')\n  ,@('6+34=', '68+26=')\n  ,@('12+34=', '9+6=')\n  ,@('60+23=', '60-0=')\n  ,@('23+68=', '48-32=')\n  ,@('28+27=', '39-28=')\n  ,@('51+19=', '87-48=')\n  ,@('58+9=', '54+8=')\n  ,@('4+47=', '20-11=')\n  ,@('46+14=', '41+22=')\n  ,@('15-11=', '54+11=')\n  ,@('34-5=', '50-9=')\n  ,@('42+11=', '7+83=')\n  ,@('10+61=', '99-90=')\n  ,@('98-76=', '20+2=')\n  ,@('95-78=', '32+22=')\n  ,@('38-19=', '95-55=')\n  ,@('68-60=', '7+26=')\n  ,@('81-64=', '34-6=')\n  ,@('26+9=', '26+12=')\n  ,@('71-49=', '90-7=')\n  ,@('17+26=', '11+9=')\n  ,@('44+14=', '97-31=')\n  ,@('82-3=', '57+41=')\n  ,@('53+38=', '59+4=')\n  ,@('27-27=', '2+78=')\n  ,@('45+51=', '47+35=')\n  ,@('31+4=', '42+18=')\n  ,@('48-26=', '54+10=')\n  ,@('62-13=', '70-54=')\n  ,@('67-16=', '93-55=')\n  ,@('10-3=', '57-12=')\n  ,@('51+33=', '16+59=')\n  ,@('91-34=', '70+10=')\n  ,@('30+69=', '88-7=')\n  ,@('48+3=', '54-21=')\n  ,@('73-41=', '62+23=')\n  ,@('6+86=', '35+41=')\n  ,@('6+48=', '1+17=')\n  ,@('79-78=', '64-28=')\n  ,@('54+29=', '74+8=')\n  ,@('25+14=', '35+52=')\n  ,@('7+70=', '7+34=')\n  ,@('55-22=', '42-37=')\n  ,@('27+37=', '24+3=')\n  ,@('70-18=', '87-77=')\n  ,@('87-76=', '27+38=')\n  ,@('70+22=', '64+10=')\n  ,@('93-51=', '37-12=')\n  ,@('24+38=', '23+11=')\n  ,@('48-13=', '97-54=')\n  ,@('10+33=', '21+11=')\n  ,@('87-3=', '9+44=')\n  ,@('92-90=', '94-67=')\n  ,@('24+4=', '43+36=')\n  ,@('23-1=', '31+46=')\n  ,@('88-60=', '52+18=')\n  ,@('64-18=', '14+26=')\n  ,@('80+5=', '6+93=')\n  ,@('9+16=', '19+80=')\n  ,@('90-53=', '24+17=')\n  ,@('98-90=', '2+30=')\n  ,@('88-17=', '41-22=')\n  ,@('14+17=', '15+31=')\n  ,@('0+83=', '51-31=')\n  ,@('29+36=', '73-69=')\n  ,@('79-48=', '24-17=')\n  ,@('63-9=', '34-17=')\n  ,@('97-66=', '34+61=')\n  ,@('92-60=', '97-38=')\n  ,@('23+26=', '72+1=')\n  ,@('2+72=', '52-1=')\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $result = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 0, $false, $newText, 2)\n  if (-not $result) {\n    Write-Output \"WARNING: replace failed for $oldText -> $newText\"\n  }\n}\n\nWrite-Output 'done'"}
